$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.739.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.27"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7899"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3150"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07009"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08050"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7582"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.900.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.296"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.786.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.920"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007668"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.147.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.138"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1641"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.282"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.047"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.378"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.389"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05686"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.065"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.590"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01902"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.770"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.824"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8396"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.024.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.898"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.851"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.437"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.054.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
